# "maxico special character done"
# Update the "Envio Semana 07" file-path references to "Envio Semana 09"
# on the Settings sheet, add a new (empty, bold-styled) cell B3, update
# the two working dates on the Constants sheet, and leave the UI focused
# on the Constants sheet/cell C24 (matching the saved view state).

$wb = $excel.ActiveWorkbook

$settings  = $wb.Worksheets.Item("Settings")
$constants = $wb.Worksheets.Item("Constants")

# ---------------------------------------------------------------------
# Settings sheet: bump the week-07 folder paths to week-09
# ---------------------------------------------------------------------
$baseDatos    = "/Planeacion/0.Envios TS/2022/09 Envío Semana 09/MÉXICO/Base de Datos"
$consolidado  = "/Planeacion/0.Envios TS/2022/09 Envio Semana 09/MÉXICO/Base de Datos/Consolidado"
$vip          = "/Planeacion/0.Envios TS/2022/09 Envío Semana 09/MÉXICO/VIP"
$tradicional  = "/Planeacion/0.Envios TS/2022/09 Envío Semana 09/MÉXICO/Tradicional"
$miCine       = "/Planeacion/0.Envios TS/2022/09 Envío Semana 09/MÉXICO/Mi Cine"
$atmosfera    = "/Planeacion/0.Envios TS/2022/09 Envío Semana 09/MÉXICO/Atmosfera"
$exportadas   = "/Planeacion/0.Envios TS/2022/09 Envío Semana 09/MÉXICO/Base de Datos/Exportadas"

$settings.Range("B2").Value  = $baseDatos
$settings.Range("B6").Value  = $baseDatos
$settings.Range("B8").Value  = $consolidado
$settings.Range("B10").Value = $vip
$settings.Range("B12").Value = $tradicional
$settings.Range("B14").Value = $miCine
$settings.Range("B16").Value = $vip
$settings.Range("B18").Value = $atmosfera
$settings.Range("B20").Value = $baseDatos
$settings.Range("B23").Value = $exportadas
$settings.Range("B26").Value = $exportadas
$settings.Range("B29").Value = $exportadas
$settings.Range("B33").Value = $exportadas

# New empty cell B3, bold-styled (creates a new bold font + cell style)
$settings.Range("B3").Font.Bold = $true

# ---------------------------------------------------------------------
# Constants sheet: move the working week forward (dates now serials)
# ---------------------------------------------------------------------
$constants.Range("B24").Value = 44613
$constants.Range("B25").Value = 44619

# ---------------------------------------------------------------------
# Final UI state: Constants tab active, selection on C24; Settings tab
# keeps a selection on B8 and scrolled near row 31.
# ---------------------------------------------------------------------
$settings.Activate()
$settings.Range("B8").Select() | Out-Null

$constants.Activate()
$constants.Range("C24").Select() | Out-Null
